$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.397.47'
$ws.Range('E2').Value = '  -1.59%  '
$ws.Range('D3').Value = '2.574.83'
$ws.Range('E3').Value = '  -3.21%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.38'
$ws.Range('E5').Value = '  -4.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '149.04'
$ws.Range('E6').Value = '  -0.72%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -1.06%  '
$ws.Range('E9').Value = '  -0.46%  '
$ws.Range('E10').Value = '  +1.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.382'
$ws.Range('E11').Value = '  -2.17%  '
$ws.Range('E12').Value = '  -0.79%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.45'
$ws.Range('E13').Value = '  -1.70%  '
$ws.Range('D14').Value = '3.036.68'
$ws.Range('E14').Value = '  -3.22%  '
$ws.Range('D15').Value = '63.242.35'
$ws.Range('E15').Value = '  -1.59%  '
$ws.Range('E16').Value = '  +3.73%  '
$ws.Range('D17').Value = '2.555.75'
$ws.Range('E17').Value = '  -3.70%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.17'
$ws.Range('E18').Value = '  +0.59%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.68'
$ws.Range('E19').Value = '  +0.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '345.21'
$ws.Range('E20').Value = '  -0.62%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.82'
$ws.Range('E21').Value = '  -1.75%  '
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.83'
$ws.Range('E23').Value = '  +0.19%  '
$ws.Range('E24').Value = '  -5.30%  '
$ws.Range('E25').Value = '  -3.83%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.08'
$ws.Range('E26').Value = '  -3.55%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '553.75'
$ws.Range('E27').Value = '  -1.54%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.162'
$ws.Range('E28').Value = '  +0.25%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.01'
$ws.Range('E29').Value = '  -3.17%  '
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('E31').Value = '  -2.82%  '
$ws.Range('D32').Value = '0.0₃0857'
$ws.Range('E32').Value = '  -0.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.75'
$ws.Range('E33').Value = '  -1.64%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.18'
$ws.Range('E34').Value = '  -2.50%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '165.22'
$ws.Range('E35').Value = '  -2.18%  '
$ws.Range('E36').Value = '  +1.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.39'
$ws.Range('E38').Value = '  -0.02%  '
$ws.Range('E39').Value = '  -3.63%  '
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '165.02'
$ws.Range('E41').Value = '  -1.59%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.73'
$ws.Range('E42').Value = '  -1.65%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.96'
$ws.Range('E43').Value = '  +2.47%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.74'
$ws.Range('E44').Value = '  +2.98%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0590'
$ws.Range('E45').Value = '  +2.32%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.626'
$ws.Range('E46').Value = '  -0.76%  '
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.04'
$ws.Range('E47').Value = '  +1.45%  '
$ws.Range('E48').Value = '  +0.81%  '
$ws.Range('E49').Value = '  -0.69%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.96'
$ws.Range('E50').Value = '  -0.66%  '
$ws.Range('D51').Value = '0.0₆0234'
$ws.Range('E51').Value = '  +16.69%  '
